$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use the same date style/format already used by the other log rows (copy from A42)
$dateStyle = $ws.Range("A42")

# New log entry (row 43) describing the Amadeus API failure
$dateStyle.Copy()
$ws.Range("A43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(43, 1).Value = 45229
$ws.Cells.Item(43, 3).Value = "Amadeus API does not work, because something is wrong with the API key, I looked at all the documentation, and I'm trying to implement some other API so hopefully solve my issues. I tried to debug it, but it kept giving a token not found error."
$ws.Cells.Item(43, 2).Value = 2

# Header rename: "time spend" -> "time spent"
$ws.Range("B1").Value = "time spent"

# New log entries (rows 45-46) describing the new flight API work
$dateStyle.Copy()
$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(45, 1).Value = 45232
$ws.Cells.Item(45, 2).Value = 3
$ws.Cells.Item(45, 3).Value = "implemented a new flight API from the aviation stack, however, I'm still getting the same errors as before where the im getting HTTP 404 not found"

$dateStyle.Copy()
$ws.Range("A46").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(46, 1).Value = 45235
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 3).Value = "reading and figured the issue could be my access point or my api key so I tried to get a new api key and that didn’t help, so I'm going to try a new access point"
$excel.CutCopyMode = 0

# Update the view to match the new scroll position / selection
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Application.ActiveWindow.Zoom = 105
$ws.Range("C46").Select()
